# Correction in SA algorithm and 746 logs:
# Replace Fitness (column C) values so that, per run, the fitness value
# reported for generations 0-24 is flattened to 8068 and for generations
# 25-250 is flattened to 7573.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row  # xlUp from bottom

for ($r = 2; $r -le $lastRow; $r++) {
    $generation = $ws.Cells.Item($r, 2).Value2

    if ($generation -le 24) {
        $ws.Cells.Item($r, 3).Value = 8068
    } else {
        $ws.Cells.Item($r, 3).Value = 7573
    }
}
